$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.60"
$ws.Range("D3").Value = "'30.96"
$ws.Range("E3").Value = "'1.58%"
$ws.Range("D4").Value = "'4.920"
$ws.Range("E4").Value = "'-0.62%"
$ws.Range("D5").Value = "'0.07318"
$ws.Range("E5").Value = "'1.87%"
$ws.Range("D6").Value = "'2.409"
$ws.Range("E6").Value = "'32.00%"
$ws.Range("D7").Value = "'7.730"
$ws.Range("E7").Value = "'0.54%"
$ws.Range("D8").Value = "'3.723"
$ws.Range("E8").Value = "'-0.92%"
$ws.Range("D9").Value = "'0.9026"
$ws.Range("D10").Value = "'0.09273"
$ws.Range("E10").Value = "'20.30%"
$ws.Range("E11").Value = "'1.77%"
$ws.Range("D12").Value = "'0.08191"
$ws.Range("E12").Value = "'2.73%"
$ws.Range("D13").Value = "'0.03127"
$ws.Range("E13").Value = "'2.82%"
$ws.Range("D14").Value = "'0.09934"
$ws.Range("E14").Value = "'-0.71%"
$ws.Range("D15").Value = "'0.001500"
$ws.Range("E15").Value = "'0.46%"
$ws.Range("D16").Value = "'0.005727"
$ws.Range("E16").Value = "'-2.06%"
$ws.Range("D17").Value = "'3.496"
$ws.Range("E17").Value = "'1.09%"
$ws.Range("D18").Value = "'2.060"
$ws.Range("E18").Value = "'-0.91%"
$ws.Range("E19").Value = "'0.30%"
$ws.Range("D20").Value = "'0.1330"
$ws.Range("E20").Value = "'4.02%"
$ws.Range("D21").Value = "'4.220"
$ws.Range("E21").Value = "'4.66%"
$ws.Range("D22").Value = "'0.2100"
$ws.Range("D23").Value = "'0.04485"
$ws.Range("E23").Value = "'-0.52%"
$ws.Range("D24").Value = "'0.001210"
$ws.Range("E24").Value = "'-0.41%"
$ws.Range("D25").Value = "'0.004159"
$ws.Range("E39").Value = "'0.52%"
$ws.Range("D40").Value = "'0.04439"
$ws.Range("E40").Value = "'1.83%"
$ws.Range("D41").Value = "'0.007357"
$ws.Range("E41").Value = "'0.90%"
$ws.Range("D42").Value = "'0.009519"
$ws.Range("E42").Value = "'-3.69%"
$ws.Range("D43").Value = "'0.1325"
$ws.Range("E43").Value = "'1.81%"
$ws.Range("D44").Value = "'0.002220"
$ws.Range("E44").Value = "'7.11%"
$ws.Range("D45").Value = "'0.008978"
$ws.Range("E45").Value = "'-5.66%"
$ws.Range("D46").Value = "'0.00006101"
$ws.Range("E46").Value = "'0.86%"
$ws.Range("E47").Value = "'-0.14%"
$ws.Range("D48").Value = "'2.575"
$ws.Range("E48").Value = "'11.50%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.14%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.14%"
